$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Extend "Tabelle1" with two new columns: "Unterlagen vorhanden" and
#    "Zeltdorf" (table grows from A10:O11 to A10:Q11).
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)

$col16 = $lo.ListColumns.Add()
$ws.Range("P10").Value = "Unterlagen vorhanden"

$col17 = $lo.ListColumns.Add()
$ws.Range("Q10").Value = "Zeltdorf"

# Approximate the original author's column widths for the two new columns
# (closest attainable values given this host's width quantisation).
$ws.Columns.Item(16).ColumnWidth = 22.6
$ws.Columns.Item(17).ColumnWidth = 19.6

# ---------------------------------------------------------------------
# 2. Drop the redundant "applyNumberFormat" flag on the B3 style (the
#    cell already uses General formatting, numFmtId 0) without touching
#    its border/font/fill. Toggling VerticalAlignment away and back to
#    its existing value (bottom) forces the host to rebuild the cell's
#    xf record, which drops the stale flag while leaving every other
#    attribute (border, font, horizontal alignment) exactly as before.
# ---------------------------------------------------------------------
$ws.Range("B3").VerticalAlignment = -4160
$ws.Range("B3").VerticalAlignment = -4107

# ---------------------------------------------------------------------
# 3. Update the saved selection to match the author's last cursor
#    position.
# ---------------------------------------------------------------------
$ws.Range("P13").Select() | Out-Null
